$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.709.28"
$ws.Range("E2").Value = "  -0.95%  "
$ws.Range("D3").Value = "2.059.60"
$ws.Range("E3").Value = "  +0.30%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.18"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.99%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.666"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.93%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "55.38"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -5.25%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "59.73"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.99%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.366"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.56%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0753"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.78%  "
$ws.Range("E12").Value = "  -3.07%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.940"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.13%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.85"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.71%  "
$ws.Range("D15").Value = "2.360.22"
$ws.Range("E15").Value = "  +0.56%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.47"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.59%  "
$ws.Range("D17").Value = "2.067.55"
$ws.Range("E17").Value = "  -0.06%  "
$ws.Range("D18").Value = "36.652.48"
$ws.Range("E18").Value = "  -0.98%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.19"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -5.62%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "72.21"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.36%  "
$ws.Range("D21").Value = "0.0₃0867"
$ws.Range("E21").Value = "  -2.29%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "238.64"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.28"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.62%  "
$ws.Range("E24").Value = "  -0.05%  "
$ws.Range("E25").Value = "  -2.53%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.16"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.37%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.34"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.88%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "165.22"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.18%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "20.22"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.90%  "
$ws.Range("E30").Value = "  -1.64%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.13"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -8.32%  "
$ws.Range("E32").Value = "  +6.24%  "
$ws.Range("E33").Value = "  -4.21%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0600"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.48%  "
$ws.Range("E35").Value = "  -0.02%  "
$ws.Range("B36").Value = "WEMIXToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.83"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.95%  "
$ws.Range("B37").Value = "Kaspa"
$ws.Range("C37").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0850"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.17%  "
$ws.Range("E38").Value = "  -1.39%  "
$ws.Range("E39").Value = "  -4.56%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.01"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.50%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.92"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.73%  "
$ws.Range("E43").Value = "  -3.85%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "94.97"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.98%  "
$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.74"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +14.54%  "
$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").Value = "1.412.46"
$ws.Range("E46").Value = "  +8.49%  "
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0906"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -7.11%  "
$ws.Range("E48").Value = "  -5.31%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.91"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.26%  "
$ws.Range("E50").Value = "  -4.14%  "
$ws.Range("D51").Value = "2.249.35"
$ws.Range("E51").Value = "  +0.89%  "
